# AFDP-7308 Combine Transcribe and OCR processing into a single media
# processing module.
#
# The "drools-ocr-rules" RuleSet previously imported the OCR-only model
# class and bound its RuleTable variable to it. Now that OCR and
# Transcribe have been folded into one media-processing module, the
# import and the RuleTable's variable declaration must point at the new
# MediaEngine model instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 (Import section): the OCR model import becomes the MediaEngine
# model import.
$ws.Range("D4").Value = "com.armedia.acm.services.mediaengine.model.MediaEngine"

# Row 17 (RuleTable OCR Rules header): the `$ocr` variable is now typed
# as MediaEngine rather than OCR.
$ws.Range("C17").Value = "`$ocr: MediaEngine"

# Leave the cursor where the edit was made, same as the author did.
$ws.Range("C17").Select()
